$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "CC"
$ws.Range("B3").Value = 1108453116
$ws.Range("C3").Value = "HP"
$ws.Range("D3").Value = 28
$ws.Range("E3").Value = "MAYO"
$ws.Range("F3").Value = 2024

$ws.Range("G2").Select()
